# Apply updated 'want-to-go' counts (column F) as scraped at commit 456a3b4
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 12   # F2: 11 -> 12
$ws.Cells.Item(3, 6).Value = 3   # F3: 2 -> 3
$ws.Cells.Item(4, 6).Value = 13255   # F4: 13245 -> 13255
$ws.Cells.Item(6, 6).Value = 2   # F6: 0 -> 2
$ws.Cells.Item(7, 6).Value = 216   # F7: 215 -> 216
$ws.Cells.Item(8, 6).Value = 113   # F8: 112 -> 113
$ws.Cells.Item(9, 6).Value = 112   # F9: 111 -> 112
$ws.Cells.Item(10, 6).Value = 62   # F10: 61 -> 62
$ws.Cells.Item(12, 6).Value = 25   # F12: 24 -> 25
$ws.Cells.Item(13, 6).Value = 13214   # F13: 13198 -> 13214
$ws.Cells.Item(14, 6).Value = 326   # F14: 325 -> 326
$ws.Cells.Item(15, 6).Value = 577   # F15: 576 -> 577
$ws.Cells.Item(16, 6).Value = 8848   # F16: 8841 -> 8848
$ws.Cells.Item(17, 6).Value = 7922   # F17: 7918 -> 7922
$ws.Cells.Item(18, 6).Value = 231   # F18: 229 -> 231
$ws.Cells.Item(21, 6).Value = 143   # F21: 142 -> 143
$ws.Cells.Item(22, 6).Value = 4   # F22: 3 -> 4
$ws.Cells.Item(23, 6).Value = 13   # F23: 12 -> 13
$ws.Cells.Item(24, 6).Value = 1003   # F24: 1001 -> 1003
$ws.Cells.Item(27, 6).Value = 389   # F27: 388 -> 389
$ws.Cells.Item(29, 6).Value = 92   # F29: 91 -> 92
$ws.Cells.Item(30, 6).Value = 356   # F30: 353 -> 356
$ws.Cells.Item(32, 6).Value = 5221   # F32: 5222 -> 5221

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 12   # F3: 11 -> 12
$ws.Cells.Item(4, 6).Value = 3   # F4: 2 -> 3
$ws.Cells.Item(5, 6).Value = 13255   # F5: 13245 -> 13255
$ws.Cells.Item(7, 6).Value = 2   # F7: 0 -> 2
$ws.Cells.Item(8, 6).Value = 216   # F8: 215 -> 216
$ws.Cells.Item(9, 6).Value = 113   # F9: 112 -> 113
$ws.Cells.Item(10, 6).Value = 112   # F10: 111 -> 112
$ws.Cells.Item(11, 6).Value = 62   # F11: 61 -> 62
$ws.Cells.Item(13, 6).Value = 25   # F13: 24 -> 25
$ws.Cells.Item(14, 6).Value = 13214   # F14: 13198 -> 13214
$ws.Cells.Item(15, 6).Value = 326   # F15: 325 -> 326
$ws.Cells.Item(16, 6).Value = 577   # F16: 576 -> 577
$ws.Cells.Item(17, 6).Value = 8848   # F17: 8841 -> 8848
$ws.Cells.Item(18, 6).Value = 7922   # F18: 7918 -> 7922
$ws.Cells.Item(19, 6).Value = 231   # F19: 229 -> 231
$ws.Cells.Item(22, 6).Value = 143   # F22: 142 -> 143
$ws.Cells.Item(23, 6).Value = 4   # F23: 3 -> 4
$ws.Cells.Item(24, 6).Value = 13   # F24: 12 -> 13
$ws.Cells.Item(25, 6).Value = 1003   # F25: 1001 -> 1003
$ws.Cells.Item(30, 6).Value = 389   # F30: 388 -> 389
$ws.Cells.Item(32, 6).Value = 92   # F32: 91 -> 92
$ws.Cells.Item(33, 6).Value = 356   # F33: 353 -> 356
$ws.Cells.Item(35, 6).Value = 5221   # F35: 5222 -> 5221
